$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corrections (Relevance Markers correction for Appenzeller-Herzog 2019 - van Dis 2020)
$ws.Range("C3").Value = 0.7916666666666666
$ws.Range("H3").Value = 0.661262050832603
$ws.Range("I3").Value = 0.07120356056144483
$ws.Range("J3").Value = 0.6944444444444444
$ws.Range("K3").Value = 681.3194444444445

$ws.Range("Q3").Value = 53
$ws.Range("R3").Value = 93
$ws.Range("S3").Value = 226
$ws.Range("T3").Value = 681
$ws.Range("U3").Value = 1579
$ws.Range("V3").Value = 9003
$ws.Range("W3").Value = 8963
$ws.Range("X3").Value = 8830
$ws.Range("Y3").Value = 8375
$ws.Range("Z3").Value = 7477

$ws.Range("AF3").Value = 0.994148
$ws.Range("AG3").Value = 0.989731
$ws.Range("AH3").Value = 0.975044
$ws.Range("AI3").Value = 0.924801
$ws.Range("AJ3").Value = 0.82564
